$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "0.9991", "4.610") keep their exact literal formatting instead of
# being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.537.80'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '1.738.17'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').Value = '0.9991'
$ws.Range('D5').Value = '247.04'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').Value = '0.9993'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').Value = '0.4902'
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('D8').Value = '0.2667'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '0.06306'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('D10').Value = '1.733.22'
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('D11').Value = '0.07037'
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = '15.73'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = '4.610'
$ws.Range('E13').Value = '  +1.88%  '
$ws.Range('D14').Value = '0.6122'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '0.9993'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '0.000007413'
$ws.Range('E17').Value = '  +7.48%  '
$ws.Range('D18').Value = '26.523.71'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = '0.9993'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '11.53'
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('D21').Value = '1.952.56'
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('D22').Value = '4.578'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').Value = '8.719'
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('D24').Value = '5.244'
$ws.Range('E24').Value = '  -1.63%  '
$ws.Range('D25').Value = '141.11'
$ws.Range('E25').Value = '  +3.88%  '
$ws.Range('D26').Value = '15.46'
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').Value = '1.415'
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('D29').Value = '108.01'
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('D30').Value = '4.038'
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('D31').Value = '0.08060'
$ws.Range('E31').Value = '  +1.40%  '
$ws.Range('D32').Value = '3.721'
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('D33').Value = '0.04585'
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '2.608'
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.010'
$ws.Range('E35').Value = '  +1.92%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.6370'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '0.8957'
$ws.Range('E37').Value = '  -4.10%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '2.018'
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.400'
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').Value = '1.005'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.01506'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '102.56'
$ws.Range('E42').Value = '  -7.42%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '5.401'
$ws.Range('E43').Value = '  -5.11%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.3903'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '6.895'
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '0.1186'
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.05394'
$ws.Range('E47').Value = '  +1.10%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').Value = '30.55'
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '7.789'
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '1.267'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '51.86'
$ws.Range('E51').Value = '  +0.84%  '
